# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rebuild the worker/period detail table (rows 16-29) on Hoja1: the two
# workers are now grouped together (all of ILMER's periods, then all of
# IVANNA's periods) instead of interleaved, and sorted with the most
# recent period (2109) first down to the oldest (2103). "Valor Mora" for
# period 2109 is 23408 for both workers; all other periods keep 35112.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$tipoDoc = "CC"
$salario = 908526

$trabajadores = @(
    @{ Doc = "12598379";   Nombre = "ILMER IVAN PASSO PUELLO" },
    @{ Doc = "1007254953"; Nombre = "IVANNA PAOLA PASSO CORREA" }
)

$periodos = @("2109", "2108", "2107", "2106", "2105", "2104", "2103")

$row = 16
foreach ($trabajador in $trabajadores) {
    foreach ($periodo in $periodos) {
        if ($periodo -eq "2109") {
            $valorMora = 23408
        } else {
            $valorMora = 35112
        }

        $ws.Cells.Item($row, 2).Value = $tipoDoc
        $ws.Cells.Item($row, 3).Value = $trabajador.Doc
        $ws.Cells.Item($row, 4).Value = $trabajador.Nombre
        $ws.Cells.Item($row, 5).Value = $periodo
        $ws.Cells.Item($row, 6).Value = $valorMora
        $ws.Cells.Item($row, 7).Value = $salario

        $row = $row + 1
    }
}

$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()
